$wb = $excel.ActiveWorkbook
Write-Output ($excel.ActiveWindow | Get-Member | Where-Object {$_.Name -like "*Tab*"})
